$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H121").Value = 939.5
$ws.Range("J121").Value = 939.5
$ws.Range("L121").Value = 2818.5
$ws.Range("N121").Value = -6312.5
$ws.Range("H131").Value = 3365.5
$ws.Range("I131").Value = 789
$ws.Range("J131").Value = 5942
$ws.Range("K131").Value = 2367
$ws.Range("L131").Value = 17826
$ws.Range("M131").Value = 2673
$ws.Range("N131").Value = -27906
$ws.Range("H132").Value = 10423086
$ws.Range("I132").Value = 12202490
$ws.Range("J132").Value = 858
$ws.Range("K132").Value = 36607470
$ws.Range("L132").Value = 2574
$ws.Range("M132").Value = -36604940
$ws.Range("N132").Value = -7634
$ws.Range("H137").Value = 3315.12
$ws.Range("I137").Value = 3098.611
$ws.Range("J137").Value = 3871.8572
$ws.Range("K137").Value = 9295.832999999999
$ws.Range("L137").Value = 11615.5716
$ws.Range("M137").Value = -6745.832999999999
$ws.Range("N137").Value = -16715.5716
$ws.Range("H141").Value = 3065.1428
$ws.Range("I141").Value = 1203.7059
$ws.Range("K141").Value = 3611.1177
$ws.Range("M141").Value = 1568.8823

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3835640.5
$ws.Range("I32").Value = 3671.1807
$ws.Range("K32").Value = 3671.1807
$ws.Range("M32").Value = -3384.1807
$ws.Range("H132").Value = 866920.75
$ws.Range("I132").Value = 1735.66
$ws.Range("J132").Value = 3270212.8
$ws.Range("K132").Value = 5206.98
$ws.Range("L132").Value = 9810638.399999999
$ws.Range("M132").Value = -2676.98
$ws.Range("N132").Value = -9815698.399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M98").ClearContents()
$ws.Range("H98").Value = 0
$ws.Range("I98").Value = 0
$ws.Range("K98").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H127").Value = 30000
$ws.Range("J127").Value = 30000
$ws.Range("L127").Value = 30000
$ws.Range("N127").Value = -39920

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 1264.625
$ws.Range("I97").Value = 1848
$ws.Range("J97").Value = 999.4545000000001
$ws.Range("K97").Value = 5544
$ws.Range("L97").Value = 2998.3635
$ws.Range("M97").Value = -5048
$ws.Range("N97").Value = -3990.3635
$ws.Range("H107").Value = 1592389.8
$ws.Range("I107").Value = 321.42856
$ws.Range("J107").Value = 2388424
$ws.Range("K107").Value = 964.28568
$ws.Range("L107").Value = 7165272
$ws.Range("M107").Value = 955.71432
$ws.Range("N107").Value = -7169112
$ws.Range("H120:N120").ClearContents()
$ws.Range("H121:N121").ClearContents()
$ws.Range("H122:N122").ClearContents()
$ws.Range("H123:N123").ClearContents()
$ws.Range("H124:N124").ClearContents()
$ws.Range("H125:M125").ClearContents()
$ws.Range("H126:M126").ClearContents()
$ws.Range("H127:L127").ClearContents()
$ws.Range("H128:L128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:M133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:L127").ClearContents()
$ws.Range("H128:L128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:L130").ClearContents()
$ws.Range("H131:L131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H135:N135").ClearContents()
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 59800
$ws.Range("J36").Value = 59800
$ws.Range("L36").Value = 59800
$ws.Range("N36").Value = -60924
